# Apply the authored changes:
#  1. Refresh the cached "datetimeFigureOut" date field text from
#     2022/12/10 -> 2023/10/7 on every slide layout and on the notes master.
#  2. Re-order the red arrow connector ("直接箭头连接符 12", id 13) on
#     slide 2 so it sits after the last shape (bring it to front / move
#     it to the end of the z-order / shape tree).

$p = $ppt.ActivePresentation

$oldDate = "2022/12/10"
$newDate = "2023/10/7"

function Update-DateShapes($shapes, $old, $new) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame -eq -1) {
            if ($shp.TextFrame.HasText -eq -1) {
                if ($shp.TextFrame.TextRange.Text -eq $old) {
                    $shp.TextFrame.TextRange.Text = $new
                }
            }
        }
    }
}

# 1a. Every slide layout off the (single) slide master.
$master = $p.SlideMaster
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DateShapes $layout.Shapes $oldDate $newDate
}

# 1b. Notes master.
$notesMaster = $p.NotesMaster
Update-DateShapes $notesMaster.Shapes $oldDate $newDate

# 2. Move the connector shape to the end of slide 2's z-order.
$slide = $p.Slides.Item(2)
for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
    $shp = $slide.Shapes.Item($i)
    if ($shp.Name -eq "直接箭头连接符 12") {
        $shp.ZOrder(0)  # msoBringToFront -> last element in the spTree
        break
    }
}
